# Update the "as_of_utc" (column AA) timestamp for the data rows
# on the "Главные" and "Линейные" sheets from 2025-11-08 03:03:59
# to 2025-11-08 07:03:51.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-08 07:03:51"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
